$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ISRG")

# Row 4: Inventory
$ws.Range("B4").Value = 577000000.0
$ws.Range("C4").Value = 602000000.0
$ws.Range("D4").Value = 663000000.0
$ws.Range("E4").Value = 646000000.0
$ws.Range("F4").Value = 620000000.0

# Row 14: Accounts Payable
$ws.Range("B14").Value = 103000000.0
$ws.Range("C14").Value = 82000000.0
$ws.Range("D14").Value = 116000000.0
$ws.Range("E14").Value = 109000000.0
$ws.Range("F14").Value = 133000000.0

# Row 22: Long Term Tax Liability (Deferred)
$ws.Range("B22").Value = -337000000.0
$ws.Range("C22").Value = -368000000.0
$ws.Range("D22").Value = -360000000.0
$ws.Range("E22").Value = -366000000.0
$ws.Range("F22").Value = -365000000.0
